# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# 1. New "Player Info" sheet (before "ODI Batting") with the player's bio.
# 2. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE (now just the
#    numeric match code instead of the full scorecard URL), plus a new
#    match row (#23) appended.
# 3. New "ODI Batting Extra" sheet (after "ODI Batting") with additional
#    per-match batting stats.

$wb = $excel.ActiveWorkbook
$wsBatting = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet: ODI Batting - rename the MATCH_CARD_LINK column to MATCH_CODE
# and replace each URL value with just the trailing MatchCode number.
# ---------------------------------------------------------------------
$wsBatting.Cells.Item(1,4).Value = "MATCH_CODE"

$lastRow = $wsBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $link = $wsBatting.Cells.Item($r,4).Value2
    $code = $link -replace '^.*MatchCode=', ''
    $wsBatting.Cells.Item($r,4).Value = "'" + $code
}

# Append new match row (match #23).
$newRow = $lastRow + 1
$rowData = @("23","23","31/03/2023","4745","1st","New Zealand","Seddon Park","run out","57","64")
for ($i = 0; $i -lt $rowData.Count; $i++) {
    $wsBatting.Cells.Item($newRow, $i + 1).Value = "'" + $rowData[$i]
}

# ---------------------------------------------------------------------
# Sheet: Player Info (new, inserted before ODI Batting)
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Add($wsBatting)
$wsInfo.Name = "Player Info"

$infoHeaders = @("ID","NAME","BATTING_HAND","BOWL_STYLE")
for ($i = 0; $i -lt $infoHeaders.Count; $i++) {
    $wsInfo.Cells.Item(1, $i + 1).Value = $infoHeaders[$i]
}

$infoRow = @("5951","Pathum Nissanka Silva","Right Handed","Does Not Bowl | Unknown")
for ($i = 0; $i -lt $infoRow.Count; $i++) {
    $wsInfo.Cells.Item(2, $i + 1).Value = "'" + $infoRow[$i]
}

# ---------------------------------------------------------------------
# Sheet: ODI Batting Extra (new, inserted after ODI Batting)
# ---------------------------------------------------------------------
$wsExtra = $wb.Worksheets.Add($null, $wb.Worksheets.Item("ODI Batting"))
$wsExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE","BATTING_POSITION","NUM_4","NUM_6","PERCENT_RUNS_OF_TOTAL","MAN_OF_MATCH")
for ($i = 0; $i -lt $extraHeaders.Count; $i++) {
    $wsExtra.Cells.Item(1, $i + 1).Value = $extraHeaders[$i]
}

$extraRows = @(
    @("4463","3","1","0","3.57%","NO"),
    @("4464","","","","","NO"),
    @("4465","3","0","0","","NO"),
    @("4469","1","0","0","2.70%","NO"),
    @("4470","1","1","0","2.07%","NO"),
    @("4471","","","","","NO"),
    @("4521","1","10","0","25.00%","NO"),
    @("4523","","","","","NO"),
    @("4527","1","6","0","21.65%","NO"),
    @("4594","","","","","NO"),
    @("4597","2","3","0","6.36%","NO"),
    @("4600","","","","","NO"),
    @("4601","2","2","0","5.04%","NO"),
    @("4603","2","0","0","1.25%","NO"),
    @("4671","","","","","NO"),
    @("4674","1","0","0","30.00%","NO"),
    @("4675","","","","","NO"),
    @("4687","","","","","NO"),
    @("4735","1","2","0","11.84%","NO"),
    @("4745","","","","","NO")
)

$rowNum = 2
foreach ($row in $extraRows) {
    $wsExtra.Cells.Item($rowNum, 1).Value = "'" + $row[0]
    if ($row[1] -ne "") {
        $wsExtra.Cells.Item($rowNum, 2).Value = [int]$row[1]
    }
    if ($row[2] -ne "") {
        $wsExtra.Cells.Item($rowNum, 3).Value = "'" + $row[2]
    }
    if ($row[3] -ne "") {
        $wsExtra.Cells.Item($rowNum, 4).Value = "'" + $row[3]
    }
    if ($row[4] -ne "") {
        $wsExtra.Cells.Item($rowNum, 5).Value = "'" + $row[4]
    }
    $wsExtra.Cells.Item($rowNum, 6).Value = "'" + $row[5]
    $rowNum++
}

# ---------------------------------------------------------------------
# Restore the first sheet as the active tab, matching the original
# workbook's activeTab state.
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
